# Append new catalogue rows (10-22) to Sheet1, matching the target commit:
# "2 logo added everything else aligned, still there's a room for improvement"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("mewa",                "language", "mewa",      "acmom",     "400 - mew"),
    @("english",              "language", "mewa",      "englis",    "400 - mew"),
    @("arabic",                "arts",     "mombasa",   "county",    "700 - mom"),
    @("the sun also rises",   "language", "hemingway", "hemingway", "400 - hem"),
    @("sentender",             "language", "emanuela",  "mark",      "400 - ema"),
    @("english",               "language", "samuel",    "etoo",      "400 - sam"),
    @("english",               "language", "arab",      "english",   "400 - ara"),
    @("english",               "language", "ann",       "etoo",      "400 - ann"),
    @("english",               "language", "ann",       "etoo",      "400 - ann"),
    @("english",               "language", "ann",       "etoo",      "400 - ann"),
    @("english",               "language", "ann",       "etoo",      "400 - ann"),
    @("english",               "language", "ann",       "etoo",      "400 - ann"),
    @("mathe",                 "language", "engli",     "ann",       "400 - eng")
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}

Write-Host "Added rows 10-22 to" $ws.Name()
